$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 404676.44; $ws.Range("I2").Value = 606190.4399999999; $ws.Range("J2").Value = 1648.5; $ws.Range("K2").Value = 606190.4399999999; $ws.Range("L2").Value = 1648.5; $ws.Range("M2").Value = -606077.4399999999; $ws.Range("N2").Value = -1874.5
$ws.Range("H3").Value = 26362.8; $ws.Range("J3").Value = 26362.8; $ws.Range("L3").Value = 26362.8; $ws.Range("N3").Value = -26590.8
$ws.Range("H19").Value = 1291.5555; $ws.Range("J19").Value = 1289.1428; $ws.Range("L19").Value = 1289.1428; $ws.Range("N19").Value = -1639.1428
$ws.Range("H62").Value = 12237.182; $ws.Range("I62").Value = 13891.556; $ws.Range("K62").Value = 13891.556; $ws.Range("M62").Value = -13267.556
$ws.Range("H65").Value = 12237.182; $ws.Range("I65").Value = 13891.556; $ws.Range("K65").Value = 69457.78; $ws.Range("M65").Value = -66337.78
$ws.Range("H74").Value = 19423; $ws.Range("J74").Value = 27398.8; $ws.Range("L74").Value = 27398.8; $ws.Range("N74").Value = -29270.8
$ws.Range("H77").Value = 19423; $ws.Range("J77").Value = 27398.8; $ws.Range("L77").Value = 136994; $ws.Range("N77").Value = -146354
$ws.Range("H102").Value = 26362.8; $ws.Range("J102").Value = 26362.8; $ws.Range("L102").Value = 26362.8; $ws.Range("N102").Value = -32852.8
$ws.Range("H103").Value = 1696.5834; $ws.Range("I103").Value = 1395; $ws.Range("J103").Value = 1724; $ws.Range("K103").Value = 4185; $ws.Range("L103").Value = 5172; $ws.Range("M103").Value = -3599; $ws.Range("N103").Value = -6344
$ws.Range("H105").Value = 19998.75; $ws.Range("J105").Value = 19998.75; $ws.Range("L105").Value = 19998.75; $ws.Range("N105").Value = -26986.75
$ws.Range("H116").Value = 11502.412; $ws.Range("I116").Value = 4232; $ws.Range("K116").Value = 4232; $ws.Range("M116").Value = -790
$ws.Range("H138").Value = 2725.625; $ws.Range("I138").Value = 2433.7368; $ws.Range("K138").Value = 7301.2104; $ws.Range("M138").Value = -2161.2104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0; $ws.Range("I10").Value = 0; $ws.Range("J10").Value = 0; $ws.Range("K10").Value = 0; $ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents(); $ws.Range("N10").ClearContents()
$ws.Range("H32").Value = 39031.656; $ws.Range("I32").Value = 42176.76; $ws.Range("K32").Value = 42176.76; $ws.Range("M32").Value = -41889.76
$ws.Range("H45").Value = 4702.25; $ws.Range("I45").Value = 1612; $ws.Range("K45").Value = 1612; $ws.Range("M45").Value = -1235
$ws.Range("H61").Value = 0; $ws.Range("I61").Value = 0; $ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H97").Value = 905.4761999999999; $ws.Range("I97").Value = 638.5454999999999; $ws.Range("J97").Value = 1199.1; $ws.Range("K97").Value = 638.5454999999999; $ws.Range("L97").Value = 1199.1; $ws.Range("M97").Value = -142.5454999999999; $ws.Range("N97").Value = -2191.1
$ws.Range("H122").Value = 2359.5881; $ws.Range("I122").Value = 1807.5333; $ws.Range("J122").Value = 6500; $ws.Range("K122").Value = 5422.5999; $ws.Range("L122").Value = 19500; $ws.Range("M122").Value = -2972.5999; $ws.Range("N122").Value = -24400
$ws.Range("H132").Value = 20977.164; $ws.Range("I132").Value = 23959.809; $ws.Range("J132").Value = 3454.125; $ws.Range("K132").Value = 71879.427; $ws.Range("L132").Value = 10362.375; $ws.Range("M132").Value = -69349.427; $ws.Range("N132").Value = -15422.375
$ws.Range("H136").Value = 0; $ws.Range("I136").Value = 0; $ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0; $ws.Range("J61").Value = 0; $ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 4303.4165; $ws.Range("J86").Value = 4735.6; $ws.Range("L86").Value = 4735.6; $ws.Range("N86").Value = -6981.6
$ws.Range("H89").Value = 4303.4165; $ws.Range("J89").Value = 4735.6; $ws.Range("L89").Value = 23678; $ws.Range("N89").Value = -34910
$ws.Range("H94").Value = 2288.2083; $ws.Range("I94").Value = 2113.5789; $ws.Range("J94").Value = 2951.8; $ws.Range("K94").Value = 2113.5789; $ws.Range("L94").Value = 2951.8; $ws.Range("M94").Value = -1662.5789; $ws.Range("N94").Value = -3853.8
$ws.Range("H95").Value = 46997.4; $ws.Range("J95").Value = 46997.4; $ws.Range("L95").Value = 46997.4; $ws.Range("N95").Value = -52489.4
$ws.Range("H99").Value = 7108.4116; $ws.Range("I99").Value = 6492.357; $ws.Range("J99").Value = 9983.333000000001; $ws.Range("K99").Value = 6492.357; $ws.Range("L99").Value = 9983.333000000001; $ws.Range("M99").Value = -4994.357; $ws.Range("N99").Value = -12979.333
$ws.Range("H105").Value = 3663.4119; $ws.Range("I105").Value = 3579.875; $ws.Range("K105").Value = 3579.875; $ws.Range("M105").Value = -1832.875
$ws.Range("H107").Value = 2166.5557; $ws.Range("I107").Value = 1571.2858; $ws.Range("J107").Value = 4250; $ws.Range("K107").Value = 1571.2858; $ws.Range("L107").Value = 4250; $ws.Range("M107").Value = 348.7141999999999; $ws.Range("N107").Value = -8090
$ws.Range("H134").Value = 4714.3076; $ws.Range("I134").Value = 3751.4707; $ws.Range("K134").Value = 11254.4121; $ws.Range("M134").Value = -8719.4121
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 200; $ws.Range("I4").Value = 100; $ws.Range("J4").Value = 300; $ws.Range("K4").Value = 100; $ws.Range("L4").Value = 300; $ws.Range("M4").Value = 12; $ws.Range("N4").Value = -524
$ws.Range("H58").Value = 49950; $ws.Range("I58").Value = 68706.53; $ws.Range("K58").Value = 68706.53; $ws.Range("M58").Value = -68503.53
$ws.Range("H136").Value = 49950; $ws.Range("I136").Value = 68706.53; $ws.Range("K136").Value = 206119.59; $ws.Range("M136").Value = -203569.59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 380370.2; $ws.Range("I4").Value = 380370.2; $ws.Range("K4").Value = 1141110.6; $ws.Range("M4").Value = -1140998.6
$ws.Range("H109").Value = 4405.857; $ws.Range("I109").Value = 2368.2; $ws.Range("K109").Value = 7104.599999999999; $ws.Range("M109").Value = -6064.599999999999
$ws.Range("H110").Value = 8000; $ws.Range("I110").Value = 8000; $ws.Range("J110").Value = 0; $ws.Range("K110").Value = 24000; $ws.Range("L110").Value = 0; $ws.Range("M110").Value = -19910
$ws.Range("N110").ClearContents()
$ws.Range("H133").Value = 10399.4; $ws.Range("I133").Value = 1999; $ws.Range("K133").Value = 5997; $ws.Range("M133").Value = -937
$ws.Range("H134").Value = 469; $ws.Range("I134").Value = 469; $ws.Range("K134").Value = 1407; $ws.Range("M134").Value = 3663
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4; $ws.Range("I5").Value = 4; $ws.Range("K5").Value = 4; $ws.Range("M5").Value = 108
$ws.Range("H70").Value = 4959.6; $ws.Range("I70").Value = 4450; $ws.Range("K70").Value = 4450; $ws.Range("M70").Value = -4180
$ws.Range("H73").Value = 4959.6; $ws.Range("I73").Value = 4450; $ws.Range("K73").Value = 4450; $ws.Range("M73").Value = -3514
$ws.Range("H97").Value = 1385.5; $ws.Range("I97").Value = 1687.6666; $ws.Range("J97").Value = 1083.3334; $ws.Range("K97").Value = 1687.6666; $ws.Range("L97").Value = 1083.3334; $ws.Range("M97").Value = -1191.6666; $ws.Range("N97").Value = -2075.3334
$ws.Range("H102").Value = 4533.3335; $ws.Range("I102").Value = 4050.25; $ws.Range("K102").Value = 4050.25; $ws.Range("M102").Value = -2428.25
$ws.Range("H126").Value = 5347.4375; $ws.Range("I126").Value = 3957.7; $ws.Range("J126").Value = 7663.6665; $ws.Range("K126").Value = 11873.1; $ws.Range("L126").Value = 22990.9995; $ws.Range("M126").Value = -9403.099999999999; $ws.Range("N126").Value = -27930.9995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 15000; $ws.Range("I12").Value = 0; $ws.Range("J12").Value = 15000; $ws.Range("K12").Value = 0; $ws.Range("L12").Value = 15000; $ws.Range("N12").Value = -15340
$ws.Range("M12").ClearContents()
$ws.Range("H133").Value = 86185.39999999999; $ws.Range("J133").Value = 86185.39999999999; $ws.Range("L133").Value = 86185.39999999999; $ws.Range("N133").Value = -91245.39999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1989.5714; $ws.Range("I122").Value = 1989.5714; $ws.Range("K122").Value = 5968.7142; $ws.Range("M122").Value = -3518.7142
$ws.Range("H133").Value = 0; $ws.Range("J133").Value = 0; $ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
